$wb = $excel.ActiveWorkbook

# --- Sheet "renters" ---
$renters = $wb.Worksheets.Item("renters")
$renters.Range("L1").Value = "Owner Savings"

# --- Sheet "owners" ---
$owners = $wb.Worksheets.Item("owners")
$owners.Range("B1").Value = "Owner Savings"

# --- Sheet "renters": add Longitude / Latitude columns ---
$renters.Range("O1").Value = "Longitude"
$renters.Range("P1").Value = "Latitude"
$renters.Range("O2:O9").Value = 0
$renters.Range("P2:P9").Value = 0
$renters.Range("O1:P1048576").Select()

# --- Sheet "owners": add Longitude / Latitude columns ---
$owners.Range("N1").Value = "Longitude"
$owners.Range("O1").Value = "Latitude"
$owners.Range("N2:N9").Value = 0
$owners.Range("O2:O9").Value = 0
$owners.Range("F2").NumberFormat = "@"
$owners.Range("F3").Select()

# --- Sheet "forsale_stock" ---
$forsale = $wb.Worksheets.Item("forsale_stock")
$forsale.Range("M1").Value = "Longitude"
$forsale.Range("N1").Value = "Latitude"
$forsale.Range("M2:M5").Value = 0
$forsale.Range("N2:N5").Value = 0
$forsale.Range("M1:N5").Select()

# --- Sheet "forrent_stock": replace "For Sale" column with Longitude / Latitude ---
$forrent = $wb.Worksheets.Item("forrent_stock")
$forrent.Range("M1").Value = "Longitude"
$forrent.Range("N1").Value = "Latitude"
$forrent.Range("M2:M5").Value = 0
$forrent.Range("N2:N5").Value = 0
$forrent.Range("F10").Select()
